$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "2018-04-10_3D_Hubs_High_Speed_Prototype.pdf" / "High speed
# prototype 3D printing" row (row 19) - everything below shifts up one row.
$ws.Rows.Item(19).Delete()

# Remove the four unused placeholder rows (now at 20-23) that only carried an
# empty, currency-formatted D cell.
$ws.Range("A20:E23").EntireRow.Delete()

# Insert three fresh rows for the new receipts, right after the
# "X-Axis high speed motor..." row (now row 19), before the Total row.
$ws.Range("A20:A22").EntireRow.Insert()

# Give the new rows the same look as a normal data row (currency format in
# column D, no special fill/border), then drop the inherited note styling in
# column E so those cells stay completely empty.
$ws.Range("A17:D17").Copy()
$ws.Range("A20:D22").PasteSpecial(-4122)
$ws.Range("E20:E22").Clear()

$ws.Range("A20").Value = "2018-06-07_Bearings_Canada.pdf"
$ws.Range("B20").Value = "1/2"" Linear bearings for high speed system"
$ws.Range("C20").Value = "Thomas"
$ws.Range("D20").Value = 99.99

$ws.Range("A21").Value = "2018-06-07_Metal_Supermarkets.pdf"
$ws.Range("B21").Value = "1/2"" Linear rods for high speed system"
$ws.Range("C21").Value = "Thomas"
$ws.Range("D21").Value = 40.91

$ws.Range("A22").Value = "2018-06-07_3D_Hubs.pdf"
$ws.Range("B22").Value = "High speed system mounting brackets"
$ws.Range("C22").Value = "Thomas"
$ws.Range("D22").Value = 72.15

$ws.Range("D17").Select() | Out-Null
